# Weekly data refresh: insert a new price record at the top of the
# historical block (row 341) for "Feria Lagunitas de Puerto Montt - Coliflor".
# All subsequent rows (old 341..372) shift down by one to (342..373).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 341, pushing existing rows 341-372 down
# to 342-373 (this also naturally grows the sheet dimension to A1:R373).
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(341, 1).Value = 4
$ws.Cells.Item(341, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(341, 3).Value = "Los Lagos"
$ws.Cells.Item(341, 4).Value = 44769
$ws.Cells.Item(341, 5).Value = 10
$ws.Cells.Item(341, 6).Value = 100112008
$ws.Cells.Item(341, 7).Value = "Coliflor"
$ws.Cells.Item(341, 8).Value = "Sin especificar"
$ws.Cells.Item(341, 9).Value = "Segunda"
$ws.Cells.Item(341, 10).Value = 250
$ws.Cells.Item(341, 11).Value = 1500
$ws.Cells.Item(341, 12).Value = 1500
$ws.Cells.Item(341, 13).Value = 1500
$ws.Cells.Item(341, 14).Value = "$/unidad"
$ws.Cells.Item(341, 15).Value = "Región Metropolitana"
$ws.Cells.Item(341, 16).Value = 1500
$ws.Cells.Item(341, 17).Value = 1
$ws.Cells.Item(341, 18).Value = "Hortaliza"
